$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2023-12-06 07:14:54"
$ws.Range("B7").Value = 0.0006000000000000001

$ws.Range("A8").Value = "2023-12-06 07:15:00"
$ws.Range("B8").Value = 0.0004
